# Updates cryptos list: new coin (OKB) inserted at row 8, pushing every
# subsequent coin down one row (Algorand drops off the bottom of the
# top-50), plus refreshed Price / Volume(1h) figures for every coin row.
#
# Price strings that look numeric (e.g. "1.001", "0.07500") are flagged
# (IsTextPrice) and written with a leading apostrophe so Excel keeps them as
# literal text -- matching the source data, which stores every price as a
# string (including multi-dot values like "29.850.46") -- then the
# resulting quote-prefix formatting is cleared so no stray number-format
# style ends up attached to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; B = 'Bitcoin'; C = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'; D = '29.850.46'; E = '  +2.37%  '; IsTextPrice = $false },
    @{ Row = 3; B = 'Ethereum'; C = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'; D = '1.862.57'; E = '  +2.11%  '; IsTextPrice = $false },
    @{ Row = 4; B = 'TetherUSD'; C = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'; D = '1.001'; E = '  +0.17%  '; IsTextPrice = $true },
    @{ Row = 5; B = 'BNB'; C = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D = '246.09'; E = '  +1.90%  '; IsTextPrice = $true },
    @{ Row = 6; B = 'XRP'; C = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'; D = '0.6415'; E = '  +3.62%  '; IsTextPrice = $true },
    @{ Row = 7; B = 'USDC'; C = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'; D = '1.001'; E = '  +0.12%  '; IsTextPrice = $true },
    @{ Row = 8; B = 'OKB'; C = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D = '47.07'; E = '  +4.72%  '; IsTextPrice = $true },
    @{ Row = 9; B = 'Cardano'; C = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D = '0.3009'; E = '  +3.90%  '; IsTextPrice = $true },
    @{ Row = 10; B = 'Dogecoin'; C = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D = '0.07500'; E = '  +2.08%  '; IsTextPrice = $true },
    @{ Row = 11; B = 'Solana'; C = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D = '24.28'; E = '  +5.48%  '; IsTextPrice = $true },
    @{ Row = 12; B = 'TRON'; C = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D = '0.07705'; E = '  +0.52%  '; IsTextPrice = $true },
    @{ Row = 13; B = 'WrappedEther'; C = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D = '1.876.24'; E = '  +2.85%  '; IsTextPrice = $false },
    @{ Row = 14; B = 'Polkadot'; C = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D = '5.076'; E = '  +2.45%  '; IsTextPrice = $true },
    @{ Row = 15; B = 'Polygon'; C = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D = '0.6892'; E = '  +4.06%  '; IsTextPrice = $true },
    @{ Row = 16; B = 'Litecoin'; C = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D = '84.25'; E = '  +2.51%  '; IsTextPrice = $true },
    @{ Row = 17; B = 'ShibaInu'; C = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D = '0.000009403'; E = '  +5.29%  '; IsTextPrice = $true },
    @{ Row = 18; B = 'Uniswap'; C = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D = '6.094'; E = '  +4.44%  '; IsTextPrice = $true },
    @{ Row = 19; B = 'WrappedBTC'; C = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D = '29.842.72'; E = '  +2.36%  '; IsTextPrice = $false },
    @{ Row = 20; B = 'WrappedliquidstakedEther2.0'; C = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D = '2.130.94'; E = '  +2.78%  '; IsTextPrice = $false },
    @{ Row = 21; B = 'BitcoinCash'; C = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D = '241.46'; E = '  +1.39%  '; IsTextPrice = $true },
    @{ Row = 22; B = 'Avalanche'; C = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D = '12.74'; E = '  +2.57%  '; IsTextPrice = $true },
    @{ Row = 23; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '1.001'; E = '  +0.11%  '; IsTextPrice = $true },
    @{ Row = 24; B = 'Chainlink'; C = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D = '7.499'; E = '  +4.17%  '; IsTextPrice = $true },
    @{ Row = 25; B = 'BinanceUSD'; C = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D = '1.002'; E = '  +0.19%  '; IsTextPrice = $true },
    @{ Row = 26; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '159.93'; E = '  +1.21%  '; IsTextPrice = $true },
    @{ Row = 27; B = 'Stellar'; C = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D = '0.1428'; E = '  +0.74%  '; IsTextPrice = $true },
    @{ Row = 28; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '8.586'; E = '  +1.60%  '; IsTextPrice = $true },
    @{ Row = 29; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '18.05'; E = '  +2.40%  '; IsTextPrice = $true },
    @{ Row = 30; B = 'PancakeSwap'; C = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D = '1.506'; E = '  +1.44%  '; IsTextPrice = $true },
    @{ Row = 31; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.06086'; E = '  +8.87%  '; IsTextPrice = $true },
    @{ Row = 32; B = 'Toncoin'; C = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D = '1.278'; E = '  +6.17%  '; IsTextPrice = $true },
    @{ Row = 33; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '4.188'; E = '  +2.32%  '; IsTextPrice = $true },
    @{ Row = 34; B = 'InternetComputer(DFINITY)'; C = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D = '4.140'; E = '  +1.09%  '; IsTextPrice = $true },
    @{ Row = 35; B = 'LidoDAOToken'; C = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D = '1.875'; E = '  +2.28%  '; IsTextPrice = $true },
    @{ Row = 36; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '1.162'; E = '  +2.65%  '; IsTextPrice = $true },
    @{ Row = 37; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '0.7352'; E = '  +0.00%  '; IsTextPrice = $true },
    @{ Row = 38; B = 'HuobiToken'; C = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D = '2.613'; E = '  -0.32%  '; IsTextPrice = $true },
    @{ Row = 39; B = 'MXToken'; C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D = '2.867'; E = '  +0.97%  '; IsTextPrice = $true },
    @{ Row = 40; B = 'VeChain'; C = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D = '0.01805'; E = '  +2.47%  '; IsTextPrice = $true },
    @{ Row = 41; B = 'Maker'; C = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'; D = '1.222.15'; E = '  +0.25%  '; IsTextPrice = $false },
    @{ Row = 42; B = 'TrustWalletToken'; C = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D = '0.9341'; E = '  +1.96%  '; IsTextPrice = $true },
    @{ Row = 43; B = 'FraxShare'; C = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D = '6.307'; E = '  +0.07%  '; IsTextPrice = $true },
    @{ Row = 44; B = 'RocketPoolETH'; C = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'; D = '2.044.14'; E = '  +3.42%  '; IsTextPrice = $false },
    @{ Row = 45; B = 'PaxDollar'; C = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D = '1.002'; E = '  +0.27%  '; IsTextPrice = $true },
    @{ Row = 46; B = 'Quant'; C = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D = '102.20'; E = '  +0.72%  '; IsTextPrice = $true },
    @{ Row = 47; B = 'Aave'; C = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D = '66.76'; E = '  +3.10%  '; IsTextPrice = $true },
    @{ Row = 48; B = 'BabyDogeCoin'; C = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'; D = '0.00000000125'; E = '  +6.25%  '; IsTextPrice = $true },
    @{ Row = 49; B = 'Mantle'; C = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D = '0.5089'; E = '  +0.13%  '; IsTextPrice = $true },
    @{ Row = 50; B = 'EnergySwap'; C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D = '9.355'; E = '  +3.13%  '; IsTextPrice = $true },
    @{ Row = 51; B = 'TheSandbox'; C = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D = '0.4106'; E = '  +2.44%  '; IsTextPrice = $true }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C

    $priceCell = $ws.Cells.Item($r.Row, 4)
    if ($r.IsTextPrice) {
        $priceCell.Value = "'" + $r.D
        $priceCell.Style = "Normal"
    } else {
        $priceCell.Value = $r.D
    }

    $ws.Cells.Item($r.Row, 5).Value = $r.E
}
